$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Location County/City" column (E) first so column indices
# for the earlier "Parent company" column (B) stay valid.
$ws.Range("E1").EntireColumn.Delete()
$ws.Range("B1").EntireColumn.Delete()

$ws.Range("A2:I3").Select()
